# "reverse pro exp and education sections"
#
# 1) Move the whole "FORMATION ACADEMIQUE" block (heading + its 4 entries) from
#    right after "PROFIL" to right before "COMPETENCES TECHNIQUES" (i.e. right
#    after the "EXPERIENCE PROFESSIONNELLE" section). This effectively swaps the
#    order of the Education and Professional-Experience sections.
# 2) Within "COMPETENCES TECHNIQUES", swap the text of the "MLOps" / "Bases de
#    donnees" bullets, and the "Langages" / "ML/AI" bullets.

$d = $word.ActiveDocument

function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    return -1
}

function Move-ParagraphBefore($needleToMove, $needleTarget) {
    $srcIdx = Find-ParagraphIndex($needleToMove)
    $src = $d.Paragraphs.Item($srcIdx)
    $rng = $d.Range($src.Range.Start, $src.Range.End)
    $rng.Cut()

    $dstIdx = Find-ParagraphIndex($needleTarget)
    $dst = $d.Paragraphs.Item($dstIdx)
    $insertPoint = $d.Range($dst.Range.Start, $dst.Range.Start)
    $insertPoint.Paste()
}

# Move each paragraph of the FORMATION ACADEMIQUE block, in its original
# order, one at a time (cut + paste individually) so that each paragraph's
# own formatting (w:pPr, e.g. the heading's border/spacing, and the entries'
# exact-line spacing) is preserved -- moving them as one multi-paragraph
# range loses the per-paragraph formatting.
$blockHeads = @(
    "FORMATION ACADEMIQUE",
    "2017-2020 : Diplôme d’ingénieur généraliste",
    "2018-2020 : Diplôme de Master IRIV",
    "2015-2017 : Classes préparatoires aux grandes écoles",
    "2014-2015 : Baccalauréat"
)

foreach ($head in $blockHeads) {
    Move-ParagraphBefore $head "COMPETENCES TECHNIQUES"
}

# Swap the MLOps <-> Bases de donnees lines, and the Langages <-> ML/AI lines,
# using a placeholder token to avoid collisions during the swap.
$placeholder = "##SWAP_PLACEHOLDER##"

$d.Content.Find.Execute(
    "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $placeholder, 2) | Out-Null

$d.Content.Find.Execute(
    "Bases de données : SQL, MongoDB, Neo4j, Redis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", 2) | Out-Null

$d.Content.Find.Execute(
    $placeholder,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bases de données : SQL, MongoDB, Neo4j, Redis", 2) | Out-Null

$d.Content.Find.Execute(
    "Langages : python, matlab, c, c++",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $placeholder, 2) | Out-Null

$d.Content.Find.Execute(
    "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Langages : python, matlab, c, c++", 2) | Out-Null

$d.Content.Find.Execute(
    $placeholder,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", 2) | Out-Null

Write-Output "done"
